$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 3 ---
$ws.Range("F3").Value = 7283749
$ws.Range("L3").Value = "D"

# --- Add new rows 4-7: additional motor vehicle license classes for Issay Paul ---

# Row 4
$ws.Range("A4").Value = "Issay"
$ws.Range("B4").Value = "Paul"
$ws.Range("C4").Value = "Kinondoni"
$ws.Range("D4").Value = "Box 4637 Kino"
$ws.Range("E4").Value = 929839
$ws.Range("F4").Value = 7283750
$ws.Range("G4").Value = "Male"
$ws.Range("H4").Value = "01/91/1992"
$ws.Range("I4").Value = "Tanzanian"
$ws.Range("J4").Value = 92839830
$ws.Range("K4").Value = "Bussines"
$ws.Range("L4").Value = "D,C"
$ws.Range("M4").Value = 42194
$ws.Range("M4").NumberFormat = $ws.Range("M3").NumberFormat

# Row 5
$ws.Range("A5").Value = "Issay"
$ws.Range("B5").Value = "Paul"
$ws.Range("C5").Value = "Kinondoni"
$ws.Range("D5").Value = "Box 4637 Kino"
$ws.Range("E5").Value = 929839
$ws.Range("F5").Value = 7283751
$ws.Range("G5").Value = "Male"
$ws.Range("H5").Value = "01/91/1992"
$ws.Range("I5").Value = "Tanzanian"
$ws.Range("J5").Value = 92839831
$ws.Range("K5").Value = "Bussines"
$ws.Range("L5").Value = "D"
$ws.Range("M5").Value = 42194
$ws.Range("M5").NumberFormat = $ws.Range("M3").NumberFormat

# Row 6
$ws.Range("A6").Value = "Issay"
$ws.Range("B6").Value = "Paul"
$ws.Range("C6").Value = "Kinondoni"
$ws.Range("D6").Value = "Box 4637 Kino"
$ws.Range("E6").Value = 929839
$ws.Range("F6").Value = 7283752
$ws.Range("G6").Value = "Male"
$ws.Range("H6").Value = "01/91/1992"
$ws.Range("I6").Value = "Tanzanian"
$ws.Range("J6").Value = 92839832
$ws.Range("K6").Value = "Bussines"
$ws.Range("L6").Value = "D,B"
$ws.Range("M6").Value = 42194
$ws.Range("M6").NumberFormat = $ws.Range("M3").NumberFormat

# Row 7
$ws.Range("A7").Value = "Issay"
$ws.Range("B7").Value = "Paul"
$ws.Range("C7").Value = "Kinondoni"
$ws.Range("D7").Value = "Box 4637 Kino"
$ws.Range("E7").Value = 929839
$ws.Range("F7").Value = 7283753
$ws.Range("G7").Value = "Male"
$ws.Range("H7").Value = "01/91/1992"
$ws.Range("I7").Value = "Tanzanian"
$ws.Range("J7").Value = 92839833
$ws.Range("K7").Value = "Bussines"
$ws.Range("L7").Value = "D"
$ws.Range("M7").Value = 42194
$ws.Range("M7").NumberFormat = $ws.Range("M3").NumberFormat

# --- Update selection to match target view state ---
$ws.Range("F9").Select()
